$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 686; this shifts the former rows 686-727
# down to 687-728 (dimension grows from A1:D727 to A1:D728).
$ws.Rows.Item(686).Insert()

# Populate the newly inserted row 686 with its data.
# Column A holds a date-like label that must stay as plain text (matching
# the surrounding cells), not get auto-converted to a date serial number.
$cellA = $ws.Range("A686")
$cellA.NumberFormat = "@"
$cellA.Value = "2026/01/20"
$cellA.NumberFormat = "General"
$cellA.Style = "Normal"

$ws.Range("B686").Value = "火"
$ws.Range("C686").Value = 6
$ws.Range("D686").Value = 140
